# Update cryptocurrency price/volume figures per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.818.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.114.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.111.18'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.44'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.484'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.629.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.800.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.112.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '475.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.70%  '
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0941'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.971'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.91'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  -2.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.311'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '387.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.832.57'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0356'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.53%  '
